$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in row 1
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 13:45"

# Refresh the raw per-country stats (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the countries whose
# figures changed in this update. Rows below refer to the CURRENT (pre-sort)
# position of each country in the sheet.
$updates = @(
    @{ Row = 4;   B = 8752121; C = 5168; D = 5698965; E = 2823839; F = 0; G = 33;  H = 229317 },  # Estados Unidos
    @{ Row = 16;  B = 562705;  C = 5814; D = 450891;  E = 79494;   F = 0; G = 335; H = 32320 },   # Iran
    @{ Row = 37;  B = 155233;  C = 2225; D = 108334;  E = 46057;   F = 0; G = 13;  H = 842 },     # Nepal
    @{ Row = 55;  B = 91978;   C = 811;  D = 82670;   E = 8355;    F = 0; G = 8;   H = 953 },     # Bielorrusia
    @{ Row = 60;  B = 78029;   C = 3614; D = 56791;   E = 20273;   F = 0; G = 11;  H = 965 },     # Austria
    @{ Row = 95;  B = 21274;   C = 1961; D = 8018;    E = 13021;   F = 0; G = 19;  H = 235 },     # Eslovenia
    @{ Row = 101; B = 15543;   C = 18;   D = 14161;   E = 1061;    F = 0; G = 0;   H = 321 },     # Senegal
    @{ Row = 134; B = 5373;    C = 115;  D = 3478;    E = 1845;    F = 0; G = 1;   H = 50 },      # Malta
    @{ Row = 144; B = 4394;    C = 86;   D = 3302;    E = 1081;    F = 0; G = 0;   H = 11 },      # Islandia
    @{ Row = 175; B = 660;     C = 19;   D = 507;     E = 153;     F = 0; G = 0;   H = 0 }        # Gibraltar
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}

# The table is kept sorted by "Casos totales" (column B) descending; re-sort
# so the countries whose totals changed land in their new rank position.
$dataRange = $ws.Range("A4:H221")
$sortKey = $ws.Range("B4:B221")
$dataRange.Sort($sortKey, 2)

# "Islas Malvinas" and "Montserrat" are tied on Casos totales (13); the source
# ranking flips their relative order on this refresh even though neither's
# figures changed, which a plain numeric sort won't reproduce on a tie.
# Pin the two rows to their final published order explicitly.
$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 2).Value = 13
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 1

$ws.Cells.Item(217, 1).Value = "Islas Malvinas"
$ws.Cells.Item(217, 2).Value = 13
$ws.Cells.Item(217, 3).Value = 0
$ws.Cells.Item(217, 4).Value = 13
$ws.Cells.Item(217, 5).Value = 0
$ws.Cells.Item(217, 6).Value = 0
$ws.Cells.Item(217, 7).Value = 0
$ws.Cells.Item(217, 8).Value = 0
